$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Add a new "email" column (U) with header + values for existing rows
# (2-11), plus two brand-new data rows (12 & 13) that replicate the
# pattern used by rows 2-11, each with its own TC_Id / email value.
# ---------------------------------------------------------------------

# New header for column U
$ws.Cells.Item(1, 21).Value = "email"

# Existing rows 2-11 -> email = "Null"
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 21).Value = "Null"
}

# Clear any stale custom row height on rows 12 & 13 before writing data
# so the resulting rows match the default (no explicit ht attribute),
# same as rows 3-11.
$ws.Rows.Item(12).AutoFit()
$ws.Rows.Item(13).AutoFit()

# ---- Row 12 ----
$ws.Cells.Item(12, 1).Value = "SI_11"
$ws.Cells.Item(12, 2).Value = "307260624P3E"
$ws.Cells.Item(12, 3).Value = "zwshashank.agrawal@teampureplay.com"
$ws.Cells.Item(12, 4).Value = 123456
$ws.Cells.Item(12, 5).Value = "userone_p1"
$ws.Cells.Item(12, 6).Value = 123456
$ws.Cells.Item(12, 7).Value = 1000
$ws.Cells.Item(12, 8).Value = 600
$ws.Cells.Item(12, 9).Value = "NULL"
$ws.Cells.Item(12, 10).Value = "NULL"
$ws.Cells.Item(12, 11).Value = "NULL"
$ws.Cells.Item(12, 12).Value = "NULL"
$ws.Cells.Item(12, 13).Value = "8906118410781 : 1"
$ws.Cells.Item(12, 14).Value = "NULL"
$ws.Cells.Item(12, 15).Value = 45384
$ws.Cells.Item(12, 16).Value = "Maharashtra"
$ws.Cells.Item(12, 17).Value = "Pune"
$ws.Cells.Item(12, 18).Value = "Dummy"
$ws.Cells.Item(12, 19).Value = "It was a great experience"
$ws.Cells.Item(12, 20).Value = 7709577438
$ws.Cells.Item(12, 21).Value = "abc@gmail.com"

# ---- Row 13 ----
$ws.Cells.Item(13, 1).Value = "SI_12"
$ws.Cells.Item(13, 2).Value = "307260624P3E"
$ws.Cells.Item(13, 3).Value = "zwshashank.agrawal@teampureplay.com"
$ws.Cells.Item(13, 4).Value = 123456
$ws.Cells.Item(13, 5).Value = "userone_p1"
$ws.Cells.Item(13, 6).Value = 123456
$ws.Cells.Item(13, 7).Value = 1000
$ws.Cells.Item(13, 8).Value = 600
$ws.Cells.Item(13, 9).Value = "NULL"
$ws.Cells.Item(13, 10).Value = "NULL"
$ws.Cells.Item(13, 11).Value = "NULL"
$ws.Cells.Item(13, 12).Value = "NULL"
$ws.Cells.Item(13, 13).Value = "8906118410781 : 1"
$ws.Cells.Item(13, 14).Value = "NULL"
$ws.Cells.Item(13, 15).Value = 45384
$ws.Cells.Item(13, 16).Value = "Maharashtra"
$ws.Cells.Item(13, 17).Value = "Pune"
$ws.Cells.Item(13, 18).Value = "Dummy"
$ws.Cells.Item(13, 19).Value = "It was a great experience"
$ws.Cells.Item(13, 20).Value = 7709577438
$ws.Cells.Item(13, 21).Value = "Null"

# ---------------------------------------------------------------------
# Fix up cell formatting that Value-only assignment cannot reproduce:
#  - I/J/K/L columns use the shaded "NULL" style seen in rows 2-11
#  - O column needs the date number format used by rows 2-11
#  - R column needs the wrap-text style used by rows 2-11
# Copy the formats (not values) from row 11 so styles match exactly.
# ---------------------------------------------------------------------
$ws.Range("I11:L11").Copy()
$ws.Range("I12:L12").PasteSpecial(-4122)
$ws.Range("I13:L13").PasteSpecial(-4122)

$ws.Range("O11").Copy()
$ws.Range("O12").PasteSpecial(-4122)
$ws.Range("O13").PasteSpecial(-4122)

$ws.Range("R11").Copy()
$ws.Range("R12").PasteSpecial(-4122)
$ws.Range("R13").PasteSpecial(-4122)

$excel.CutCopyMode = 0
